$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Population")

# Row 3 (scenario 1) was a duplicate of row 2 (scenario 1) - remove it.
# Deleting the row shifts all subsequent rows up one position.
$ws.Rows.Item(3).Delete()

# Renumber the scenario id column (A) sequentially 1..N now that the
# duplicate row has been removed.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 1
}

# Update the selection on the Population sheet and make it the active tab.
$ws.Range("J17").Select()
$ws.Activate()
